{"js": "// Insert a new \"My name : Erifeoluwa Jamgbadi\" paragraph (preceded by a new\n// blank paragraph) right before the empty paragraph that immediately\n// precedes the \"Answers to/solutions of the single exercises in the\n// assignment\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the \"Answers to/solutions...\" paragraph; the new content goes\n// right before the blank paragraph that sits just above it.\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Answers to/solutions\") !== -1) {\n    anchorIndex = i - 1;\n    break;\n  }\n}\nif (anchorIndex < 0) {\n  throw new Error(\"Could not locate insertion anchor paragraph\");\n}\nconst target = paragraphs.items[anchorIndex];\n\n// Create a placeholder paragraph before the target; its content will be\n// replaced with the real run/proofErr structure below.\nconst placeholder = target.insertParagraph(\"PLACEHOLDER\", Word.InsertLocation.before);\nawait context.sync();\n\n// Build the OOXML for the two new paragraphs: an empty paragraph followed\n// by the \"My name : Erifeoluwa Jamgbadi\" paragraph split across three runs\n// (matching Word's grammar-check proofErr markers around \"name :\").\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p/>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">My </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t>name :</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> Erifeoluwa Jamgbadi</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst placeholderRange = placeholder.getRange(\"Content\");\nplaceholderRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Insert a new \"My name : Erifeoluwa Jamgbadi\" paragraph (preceded by a new\n# blank paragraph) right before the empty paragraph that immediately\n# precedes the \"Answers to/solutions of the single exercises in the\n# assignment\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Answers to/solutions...\" paragraph; the new content goes\n# right before the blank paragraph that sits just above it.\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Answers to/solutions*\") {\n        $anchorIndex = $i - 1\n        break\n    }\n}\nif ($anchorIndex -lt 1) {\n    throw \"Could not locate insertion anchor paragraph\"\n}\n\n$target = $d.Paragraphs.Item($anchorIndex)\n$targetRange = $target.Range\n\n# Create a placeholder paragraph before the target; its content will be\n# replaced with the real run/proofErr structure below.\n$targetRange.InsertParagraphBefore()\n\n$placeholder = $d.Paragraphs.Item($anchorIndex)\n$placeholderRange = $placeholder.Range\n\n# Build the OOXML for the two new paragraphs: an empty paragraph followed\n# by the \"My name : Erifeoluwa Jamgbadi\" paragraph split across three runs\n# (matching Word's grammar-check proofErr markers around \"name :\").\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p/><w:p><w:r><w:t xml:space=\"preserve\">My </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>name :</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> Erifeoluwa Jamgbadi</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$placeholderRange.InsertXML($xml)\n"}
